$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for each data row (2-78). Rows 2-76 are a reshuffle of the existing
# 75 data rows; rows 77-78 are brand new rows appended at the bottom.
$table = @{}
$table[2] = @(44357, 3200, 1300, 1500, 1400, 933)
$table[3] = @(44243, 3200, 1500, 2000, 1750, 1167)
$table[4] = @(44350, 3340, 1200, 1500, 1350, 900)
$table[5] = @(44208, 3400, 1300, 1500, 1400, 933)
$table[6] = @(44327, 3400, 1300, 1500, 1400, 933)
$table[7] = @(44280, 3000, 2000, 2500, 2250, 1500)
$table[8] = @(44383, 3200, 1500, 2000, 1750, 1167)
$table[9] = @(44364, 3200, 1500, 2000, 1750, 1167)
$table[10] = @(44336, 3360, 1300, 1500, 1400, 933)
$table[11] = @(44306, 3400, 2000, 2500, 2250, 1500)
$table[12] = @(44285, 3400, 2000, 2500, 2250, 1500)
$table[13] = @(44390, 3200, 1500, 2000, 1750, 1167)
$table[14] = @(44201, 2800, 1300, 1500, 1400, 933)
$table[15] = @(44425, 3200, 2000, 2500, 2250, 1500)
$table[16] = @(44266, 3600, 2000, 2500, 2250, 1500)
$table[17] = @(44166, 2800, 1300, 1500, 1400, 933)
$table[18] = @(44413, 3360, 2000, 2500, 2250, 1500)
$table[19] = @(44308, 3200, 1300, 1500, 1400, 933)
$table[20] = @(44271, 3200, 2000, 2500, 2250, 1500)
$table[21] = @(44229, 3200, 1300, 1500, 1400, 933)
$table[22] = @(44189, 2600, 1400, 1500, 1450, 967)
$table[23] = @(44231, 2800, 1300, 1500, 1400, 933)
$table[24] = @(44355, 3200, 1300, 1500, 1400, 933)
$table[25] = @(44203, 2800, 1300, 1500, 1400, 933)
$table[26] = @(44418, 3300, 2000, 2500, 2250, 1500)
$table[27] = @(44236, 3200, 1500, 2000, 1750, 1167)
$table[28] = @(44159, 2900, 1000, 1500, 1250, 833)
$table[29] = @(44245, 3200, 1500, 2000, 1750, 1167)
$table[30] = @(44187, 3100, 1400, 1500, 1450, 967)
$table[31] = @(44343, 3340, 1300, 1500, 1400, 933)
$table[32] = @(44320, 3400, 1300, 1500, 1400, 933)
$table[33] = @(44278, 3400, 2000, 2500, 2250, 1500)
$table[34] = @(44406, 3400, 2000, 2500, 2250, 1500)
$table[35] = @(44385, 3320, 1500, 2000, 1750, 1167)
$table[36] = @(44250, 3400, 1500, 2000, 1750, 1167)
$table[37] = @(44334, 3440, 1300, 1500, 1400, 933)
$table[38] = @(44299, 3400, 2000, 2500, 2250, 1500)
$table[39] = @(44252, 3600, 1500, 2000, 1750, 1167)
$table[40] = @(44371, 3300, 1500, 2000, 1750, 1167)
$table[41] = @(44292, 3400, 2000, 2500, 2250, 1500)
$table[42] = @(44259, 3400, 2000, 2500, 2250, 1500)
$table[43] = @(44392, 3320, 1500, 2000, 1750, 1167)
$table[44] = @(44224, 2800, 1300, 1500, 1400, 933)
$table[45] = @(44217, 2800, 1300, 1500, 1400, 933)
$table[46] = @(44362, 3200, 1500, 2000, 1750, 1167)
$table[47] = @(44168, 2800, 1300, 1500, 1400, 933)
$table[48] = @(44322, 3320, 1300, 1500, 1400, 933)
$table[49] = @(44196, 3200, 1400, 1500, 1450, 967)
$table[50] = @(44397, 3200, 1500, 2000, 1750, 1167)
$table[51] = @(44434, 3360, 2000, 2500, 2250, 1500)
$table[52] = @(44215, 2800, 1300, 1500, 1400, 933)
$table[53] = @(44194, 3300, 1400, 1500, 1450, 967)
$table[54] = @(44222, 2800, 1300, 1500, 1400, 933)
$table[55] = @(44420, 3400, 2000, 2500, 2250, 1500)
$table[56] = @(44264, 3600, 2000, 2500, 2250, 1500)
$table[57] = @(44348, 3360, 1300, 1500, 1400, 933)
$table[58] = @(44210, 3200, 1300, 1500, 1400, 933)
$table[59] = @(44427, 3360, 2000, 2500, 2250, 1500)
$table[60] = @(44341, 3360, 1300, 1500, 1400, 933)
$table[61] = @(44301, 3200, 2000, 2500, 2250, 1500)
$table[62] = @(44441, 3200, 2000, 2500, 2250, 1500)
$table[63] = @(44432, 3200, 2000, 2500, 2250, 1500)
$table[64] = @(44369, 3200, 1500, 2000, 1750, 1167)
$table[65] = @(44257, 3600, 2000, 2500, 2250, 1500)
$table[66] = @(44273, 3000, 2000, 2500, 2250, 1500)
$table[67] = @(44294, 3000, 2000, 2500, 2250, 1500)
$table[68] = @(44315, 3120, 1300, 1500, 1400, 933)
$table[69] = @(44446, 3200, 2000, 2500, 2250, 1500)
$table[70] = @(44411, 3300, 2000, 2500, 2250, 1500)
$table[71] = @(44313, 3200, 1300, 1500, 1400, 933)
$table[72] = @(44329, 3300, 1300, 1500, 1400, 933)
$table[73] = @(44161, 3100, 1300, 1500, 1400, 933)
$table[74] = @(44435, 6560, 2000, 2500, 2250, 1500)
$table[75] = @(44175, 3000, 1300, 1500, 1400, 933)
$table[76] = @(44376, 3200, 1500, 2000, 1750, 1167)
$table[77] = @(44238, 3200, 1500, 2000, 1750, 1167)
$table[78] = @(44399, 3320, 1500, 2000, 1750, 1167)

# Static values shared by every data row (unchanged by this edit)
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112044
$categoria = "Perejil"
$variedad = "Sin especificar"
$calidad = "Primera"
$unidad = "`$/atado 1 a 1,5 kilos"
$origen = "Provincia del Elqu" + [char]0x00ED
$kgUnidades = 1.5
$clasificacion = "Hortaliza"

for ($r = 2; $r -le 78; $r++) {
    $vals = $table[$r]

    if ($r -ge 77) {
        # Brand new rows: populate every column.
        $ws.Cells.Item($r, 1).Value = $mercadoId
        $ws.Cells.Item($r, 2).Value = $mercado
        $ws.Cells.Item($r, 3).Value = $region
        $ws.Cells.Item($r, 5).Value = $codreg
        $ws.Cells.Item($r, 6).Value = $categoriaId
        $ws.Cells.Item($r, 7).Value = $categoria
        $ws.Cells.Item($r, 8).Value = $variedad
        $ws.Cells.Item($r, 9).Value = $calidad
        $ws.Cells.Item($r, 14).Value = $unidad
        $ws.Cells.Item($r, 15).Value = $origen
        $ws.Cells.Item($r, 17).Value = $kgUnidades
        $ws.Cells.Item($r, 18).Value = $clasificacion
    }

    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 10).Value = $vals[1]
    $ws.Cells.Item($r, 11).Value = $vals[2]
    $ws.Cells.Item($r, 12).Value = $vals[3]
    $ws.Cells.Item($r, 13).Value = $vals[4]
    $ws.Cells.Item($r, 16).Value = $vals[5]
}

Write-Output "done"
